$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 10744.615
$ws.Cells.Item(86, 9).Value = 14092.5
$ws.Cells.Item(86, 10).Value = 5388
$ws.Cells.Item(86, 11).Value = 14092.5
$ws.Cells.Item(86, 12).Value = 5388
$ws.Cells.Item(86, 13).Value = -12969.5
$ws.Cells.Item(86, 14).Value = -7634
$ws.Cells.Item(89, 8).Value = 10744.615
$ws.Cells.Item(89, 9).Value = 14092.5
$ws.Cells.Item(89, 10).Value = 5388
$ws.Cells.Item(89, 11).Value = 70462.5
$ws.Cells.Item(89, 12).Value = 26940
$ws.Cells.Item(89, 13).Value = -64846.5
$ws.Cells.Item(89, 14).Value = -38172
$ws.Cells.Item(98, 8).Value = 1470.9412
$ws.Cells.Item(98, 9).Value = 1470.9412
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 1470.9412
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 27.05880000000002
$ws.Cells.Item(98, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 1897.1666
$ws.Cells.Item(107, 9).Value = 2842
$ws.Cells.Item(107, 10).Value = 1222.2858
$ws.Cells.Item(107, 11).Value = 2842
$ws.Cells.Item(107, 12).Value = 1222.2858
$ws.Cells.Item(107, 13).Value = -922
$ws.Cells.Item(107, 14).Value = -5062.2858
$ws.Cells.Item(116, 8).Value = 6452.222
$ws.Cells.Item(116, 9).Value = 5894
$ws.Cells.Item(116, 10).Value = 7150
$ws.Cells.Item(116, 11).Value = 5894
$ws.Cells.Item(116, 12).Value = 7150
$ws.Cells.Item(116, 13).Value = -2452
$ws.Cells.Item(116, 14).Value = -14034
$ws.Cells.Item(122, 8).Value = 1470.9412
$ws.Cells.Item(122, 9).Value = 1470.9412
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4412.8236
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1962.8236
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(125, 8).Value = 1391000.5
$ws.Cells.Item(125, 10).Value = 1391000.5
$ws.Cells.Item(125, 12).Value = 12519004.5
$ws.Cells.Item(125, 14).Value = -12523924.5
$ws.Cells.Item(127, 8).Value = 1198.5
$ws.Cells.Item(127, 9).Value = 397
$ws.Cells.Item(127, 10).Value = 2000
$ws.Cells.Item(127, 11).Value = 1191
$ws.Cells.Item(127, 12).Value = 6000
$ws.Cells.Item(127, 13).Value = 3769
$ws.Cells.Item(127, 14).Value = -15920
$ws.Cells.Item(132, 8).Value = 6581622.5
$ws.Cells.Item(132, 9).Value = 8198228
$ws.Cells.Item(132, 10).Value = 7427.1333
$ws.Cells.Item(132, 11).Value = 24594684
$ws.Cells.Item(132, 12).Value = 22281.3999
$ws.Cells.Item(132, 13).Value = -24592154
$ws.Cells.Item(132, 14).Value = -27341.3999
$ws.Cells.Item(133, 8).Value = 53385
$ws.Cells.Item(133, 10).Value = 53385
$ws.Cells.Item(133, 12).Value = 53385
$ws.Cells.Item(133, 14).Value = -63505
$ws.Cells.Item(135, 8).Value = 214586.66
$ws.Cells.Item(135, 9).Value = 272141.7
$ws.Cells.Item(135, 10).Value = 1633
$ws.Cells.Item(135, 11).Value = 2449275.3
$ws.Cells.Item(135, 12).Value = 14697
$ws.Cells.Item(135, 13).Value = -2446740.3
$ws.Cells.Item(135, 14).Value = -19767
$ws.Cells.Item(137, 8).Value = 2245.0715
$ws.Cells.Item(137, 9).Value = 1221
$ws.Cells.Item(137, 10).Value = 6000
$ws.Cells.Item(137, 11).Value = 3663
$ws.Cells.Item(137, 12).Value = 18000
$ws.Cells.Item(137, 13).Value = -1113
$ws.Cells.Item(137, 14).Value = -23100
$ws.Cells.Item(138, 8).Value = 6890.5713
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 6890.5713
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 20671.7139
$ws.Cells.Item(138, 13).ClearContents()
$ws.Cells.Item(138, 14).Value = -30951.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(6, 8).Value = 15001
$ws.Cells.Item(6, 9).Value = 20000
$ws.Cells.Item(6, 10).Value = 12501.5
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 12501.5
$ws.Cells.Item(6, 13).Value = -19827
$ws.Cells.Item(6, 14).Value = -12847.5
$ws.Cells.Item(9, 8).Value = 18000
$ws.Cells.Item(9, 10).Value = 18000
$ws.Cells.Item(9, 12).Value = 18000
$ws.Cells.Item(9, 14).Value = -18340
$ws.Cells.Item(20, 8).Value = 18000
$ws.Cells.Item(20, 10).Value = 18000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 14).Value = -18540
$ws.Cells.Item(37, 8).Value = 10979.2
$ws.Cells.Item(37, 9).Value = 15000
$ws.Cells.Item(37, 10).Value = 9974
$ws.Cells.Item(37, 11).Value = 15000
$ws.Cells.Item(37, 12).Value = 9974
$ws.Cells.Item(37, 13).Value = -14727
$ws.Cells.Item(37, 14).Value = -10520
$ws.Cells.Item(44, 8).Value = 24673.166
$ws.Cells.Item(44, 10).Value = 34009.75
$ws.Cells.Item(44, 12).Value = 34009.75
$ws.Cells.Item(44, 14).Value = -34985.75
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 14).ClearContents()
$ws.Cells.Item(63, 8).Value = 3614.2856
$ws.Cells.Item(63, 9).Value = 2233.3333
$ws.Cells.Item(63, 10).Value = 4166.6665
$ws.Cells.Item(63, 11).Value = 2233.3333
$ws.Cells.Item(63, 12).Value = 4166.6665
$ws.Cells.Item(63, 13).Value = -1547.3333
$ws.Cells.Item(63, 14).Value = -5538.6665
$ws.Cells.Item(66, 8).Value = 3614.2856
$ws.Cells.Item(66, 9).Value = 2233.3333
$ws.Cells.Item(66, 10).Value = 4166.6665
$ws.Cells.Item(66, 11).Value = 11166.6665
$ws.Cells.Item(66, 12).Value = 20833.3325
$ws.Cells.Item(66, 13).Value = -7734.666499999999
$ws.Cells.Item(66, 14).Value = -27697.3325
$ws.Cells.Item(74, 8).Value = 22201.766
$ws.Cells.Item(74, 9).Value = 2164.1714
$ws.Cells.Item(74, 10).Value = 66034
$ws.Cells.Item(74, 11).Value = 2164.1714
$ws.Cells.Item(74, 12).Value = 66034
$ws.Cells.Item(74, 13).Value = -1290.1714
$ws.Cells.Item(74, 14).Value = -67782
$ws.Cells.Item(77, 8).Value = 22201.766
$ws.Cells.Item(77, 9).Value = 2164.1714
$ws.Cells.Item(77, 10).Value = 66034
$ws.Cells.Item(77, 11).Value = 10820.857
$ws.Cells.Item(77, 12).Value = 330170
$ws.Cells.Item(77, 13).Value = -6452.857
$ws.Cells.Item(77, 14).Value = -338906
$ws.Cells.Item(80, 8).Value = 38332.668
$ws.Cells.Item(80, 10).Value = 38332.668
$ws.Cells.Item(80, 12).Value = 38332.668
$ws.Cells.Item(80, 14).Value = -40328.668
$ws.Cells.Item(83, 8).Value = 38332.668
$ws.Cells.Item(83, 10).Value = 38332.668
$ws.Cells.Item(83, 12).Value = 114998.004
$ws.Cells.Item(83, 14).Value = -124982.004
$ws.Cells.Item(132, 8).Value = 2498.0212
$ws.Cells.Item(132, 9).Value = 1485.9524
$ws.Cells.Item(132, 10).Value = 10999.4
$ws.Cells.Item(132, 11).Value = 4457.857199999999
$ws.Cells.Item(132, 12).Value = 32998.2
$ws.Cells.Item(132, 13).Value = -1927.857199999999
$ws.Cells.Item(132, 14).Value = -38058.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(51, 8).Value = 14350
$ws.Cells.Item(51, 10).Value = 14350
$ws.Cells.Item(51, 12).Value = 14350
$ws.Cells.Item(51, 14).Value = -15332
$ws.Cells.Item(82, 8).Value = 24703.688
$ws.Cells.Item(82, 10).Value = 44055.5
$ws.Cells.Item(82, 12).Value = 44055.5
$ws.Cells.Item(82, 14).Value = -44821.5
$ws.Cells.Item(85, 8).Value = 24703.688
$ws.Cells.Item(85, 10).Value = 44055.5
$ws.Cells.Item(85, 12).Value = 44055.5
$ws.Cells.Item(85, 14).Value = -46707.5
$ws.Cells.Item(134, 8).Value = 1442.7954
$ws.Cells.Item(134, 9).Value = 1421.1951
$ws.Cells.Item(134, 11).Value = 4263.5853
$ws.Cells.Item(134, 13).Value = -1728.5853

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 7861.125
$ws.Cells.Item(50, 10).Value = 8972.286
$ws.Cells.Item(50, 12).Value = 8972.286
$ws.Cells.Item(50, 14).Value = -10222.286
$ws.Cells.Item(51, 8).Value = 9385.200000000001
$ws.Cells.Item(51, 10).Value = 9385.200000000001
$ws.Cells.Item(51, 12).Value = 9385.200000000001
$ws.Cells.Item(51, 14).Value = -10857.2
$ws.Cells.Item(60, 8).Value = 13282.333
$ws.Cells.Item(60, 10).Value = 13282.333
$ws.Cells.Item(60, 12).Value = 13282.333
$ws.Cells.Item(60, 14).Value = -14304.333
$ws.Cells.Item(61, 8).Value = 9385.200000000001
$ws.Cells.Item(61, 10).Value = 9385.200000000001
$ws.Cells.Item(61, 12).Value = 9385.200000000001
$ws.Cells.Item(61, 14).Value = -10081.2
$ws.Cells.Item(68, 8).Value = 15020
$ws.Cells.Item(68, 10).Value = 17482
$ws.Cells.Item(68, 12).Value = 17482
$ws.Cells.Item(68, 14).Value = -18980
$ws.Cells.Item(71, 8).Value = 15020
$ws.Cells.Item(71, 10).Value = 17482
$ws.Cells.Item(71, 12).Value = 52446
$ws.Cells.Item(71, 14).Value = -59934
$ws.Cells.Item(107, 8).Value = 1000.2
$ws.Cells.Item(107, 9).Value = 1411.7142
$ws.Cells.Item(107, 10).Value = 640.125
$ws.Cells.Item(107, 11).Value = 1411.7142
$ws.Cells.Item(107, 12).Value = 640.125
$ws.Cells.Item(107, 13).Value = 508.2858000000001
$ws.Cells.Item(107, 14).Value = -4480.125
$ws.Cells.Item(109, 8).Value = 13618.333
$ws.Cells.Item(109, 10).Value = 13618.333
$ws.Cells.Item(109, 12).Value = 13618.333
$ws.Cells.Item(109, 14).Value = -15698.333
$ws.Cells.Item(118, 8).Value = 33000
$ws.Cells.Item(118, 10).Value = 33000
$ws.Cells.Item(118, 12).Value = 33000
$ws.Cells.Item(118, 14).Value = -36314
$ws.Cells.Item(132, 8).Value = 1900.0536
$ws.Cells.Item(132, 9).Value = 1605.9302
$ws.Cells.Item(132, 10).Value = 2872.923
$ws.Cells.Item(132, 11).Value = 4817.7906
$ws.Cells.Item(132, 12).Value = 8618.769
$ws.Cells.Item(132, 13).Value = -2287.7906
$ws.Cells.Item(132, 14).Value = -13678.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 9).Value = 3000
$ws.Cells.Item(122, 11).Value = 9000
$ws.Cells.Item(122, 13).Value = -6550
$ws.Cells.Item(133, 8).Value = 23055.666
$ws.Cells.Item(133, 10).Value = 24687.625
$ws.Cells.Item(133, 12).Value = 24687.625
$ws.Cells.Item(133, 14).Value = -34807.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(38, 8).Value = 4033
$ws.Cells.Item(38, 10).Value = 4033
$ws.Cells.Item(38, 12).Value = 4033
$ws.Cells.Item(38, 14).Value = -4853
$ws.Cells.Item(130, 8).Value = 52463
$ws.Cells.Item(130, 10).Value = 52463
$ws.Cells.Item(130, 12).Value = 52463
$ws.Cells.Item(130, 14).Value = -62503
$ws.Cells.Item(132, 8).Value = 4778.1113
$ws.Cells.Item(132, 9).Value = 4550.6665
$ws.Cells.Item(132, 10).Value = 5233
$ws.Cells.Item(132, 11).Value = 13651.9995
$ws.Cells.Item(132, 12).Value = 15699
$ws.Cells.Item(132, 13).Value = -11121.9995
$ws.Cells.Item(132, 14).Value = -20759

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3047.2083
$ws.Cells.Item(132, 9).Value = 3318.3845
$ws.Cells.Item(132, 10).Value = 2726.7273
$ws.Cells.Item(132, 11).Value = 9955.1535
$ws.Cells.Item(132, 12).Value = 8180.1819
$ws.Cells.Item(132, 13).Value = -7425.1535
$ws.Cells.Item(132, 14).Value = -13240.1819
